# "Generate Report for Handoff"
# The localization-status report moved from "In Translation" to
# "Ready for handoff" and a fresh handoff xliff-generation timestamp was
# stamped on every sheet (Overview + each locale tab). Excel re-flowed the
# Status/Datetime columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-13 19:14:19"
$ws.Columns.Item(5).ColumnWidth = 16.4
$ws.Columns.Item(6).ColumnWidth = 16.4

# ---- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-13 19:14:11"
$ws.Columns.Item(3).ColumnWidth = 16.4

# ---- de-de sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-13 19:14:19"
$ws.Columns.Item(3).ColumnWidth = 16.4
